# Cam_SpecialEventLobbyist.xlsx edit
# 1) Fix TransactionStatusId (col I) on a handful of existing rows: 1 -> 2
# 2) Append 26 new lobbyist-event rows (2278-2303)
# 3) Extend the "Cam_SpecialEventLobbyist" defined name to the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix TransactionStatusId (column I) for a handful of existing rows: 1 -> 2 ---
$fixRows = @(2259, 2261, 2262, 2263, 2264, 2277)
foreach ($r in $fixRows) {
    $ws.Cells.Item($r, 9).Value = 2
}

# --- Append new lobbyist data rows 2278-2303 ---
# Row 2278
$ws.Cells.Item(2278, 1).Value = 2478
$ws.Cells.Item(2278, 2).Value = 8357
$ws.Cells.Item(2278, 3).Value = "Committee Dinner"
$ws.Cells.Item(2278, 4).Value = "Bull Ring"
$ws.Cells.Item(2278, 5).Value = 43451
$ws.Cells.Item(2278, 6).Value = 235.49
$ws.Cells.Item(2278, 7).Value = "Revenue Stabilization committee and staff"
$ws.Cells.Item(2278, 8).Value = 43452.5665096412
$ws.Cells.Item(2278, 9).Value = 2
$ws.Cells.Item(2278, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2278, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2278, 8).NumberFormat = "mm-dd-yy"

# Row 2279
$ws.Cells.Item(2279, 1).Value = 2479
$ws.Cells.Item(2279, 2).Value = 8435
$ws.Cells.Item(2279, 3).Value = "Committee dinner"
$ws.Cells.Item(2279, 4).Value = "Pepper's Grill and Bar"
$ws.Cells.Item(2279, 5).Value = 43375
$ws.Cells.Item(2279, 6).Value = 734.23
$ws.Cells.Item(2279, 7).Value = "Water & Natural Resources interim committee and staff"
$ws.Cells.Item(2279, 8).Value = 43455.7501158565
$ws.Cells.Item(2279, 9).Value = 1
$ws.Cells.Item(2279, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2279, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2279, 8).NumberFormat = "mm-dd-yy"

# Row 2280
$ws.Cells.Item(2280, 1).Value = 2480
$ws.Cells.Item(2280, 2).Value = 8462
$ws.Cells.Item(2280, 3).Value = "Dinner"
$ws.Cells.Item(2280, 4).Value = "Salud de Mesilla"
$ws.Cells.Item(2280, 5).Value = 43384
$ws.Cells.Item(2280, 6).Value = 186.43
$ws.Cells.Item(2280, 7).Value = "Members, staff and guests of the Economic Development Committee."
$ws.Cells.Item(2280, 8).Value = 43468.4978069444
$ws.Cells.Item(2280, 9).Value = 2
$ws.Cells.Item(2280, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2280, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2280, 8).NumberFormat = "mm-dd-yy"

# Row 2281
$ws.Cells.Item(2281, 1).Value = 2481
$ws.Cells.Item(2281, 2).Value = 8462
$ws.Cells.Item(2281, 3).Value = "Lunch"
$ws.Cells.Item(2281, 4).Value = "State Capitol"
$ws.Cells.Item(2281, 5).Value = 43451
$ws.Cells.Item(2281, 6).Value = 369.98
$ws.Cells.Item(2281, 7).Value = "RSTP members, staff and guests."
$ws.Cells.Item(2281, 8).Value = 43468.4986045139
$ws.Cells.Item(2281, 9).Value = 2
$ws.Cells.Item(2281, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2281, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2281, 8).NumberFormat = "mm-dd-yy"

# Row 2282
$ws.Cells.Item(2282, 1).Value = 2482
$ws.Cells.Item(2282, 2).Value = 8491
$ws.Cells.Item(2282, 3).Value = "AED Legislative Luncheon"
$ws.Cells.Item(2282, 4).Value = "Artichoke Cafe"
$ws.Cells.Item(2282, 5).Value = 43382
$ws.Cells.Item(2282, 6).Value = 1023
$ws.Cells.Item(2282, 7).Value = "Selected legislators were invited as well as candidates for legislative seats and members of the AED Board of Directors. The legislators who attended were:_x000d_`nDaymon Ely_x000d_`nJimmie Hall_x000d_`nDaniel Ivey-Soto_x000d_`nAntonio Maestas_x000d_`nBill Payne_x000d_`nGregg Schmedes_x000d_`nMimi Stewart_x000d_`nJames White_x000d_`nThe candidates for legislative seats who attended were:_x000d_`nRobert Godshall_x000d_`nJohn Jones_x000d_`nDayan Hochman-Vigil"
$ws.Cells.Item(2282, 8).Value = 43473
$ws.Cells.Item(2282, 9).Value = 2
$ws.Cells.Item(2282, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2282, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2282, 8).NumberFormat = "mm-dd-yy"

# Row 2283
$ws.Cells.Item(2283, 1).Value = 2483
$ws.Cells.Item(2283, 2).Value = 8515
$ws.Cells.Item(2283, 3).Value = "Committee Lunch"
$ws.Cells.Item(2283, 4).Value = "Capitol Building"
$ws.Cells.Item(2283, 5).Value = 43405
$ws.Cells.Item(2283, 6).Value = 258.98
$ws.Cells.Item(2283, 7).Value = "Economic and Rural Development Committee and Staff"
$ws.Cells.Item(2283, 8).Value = 43474.514821794
$ws.Cells.Item(2283, 9).Value = 2
$ws.Cells.Item(2283, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2283, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2283, 8).NumberFormat = "mm-dd-yy"

# Row 2284
$ws.Cells.Item(2284, 1).Value = 2484
$ws.Cells.Item(2284, 2).Value = 8515
$ws.Cells.Item(2284, 3).Value = "Committee Lunch"
$ws.Cells.Item(2284, 4).Value = "Capitol Building"
$ws.Cells.Item(2284, 5).Value = 43433
$ws.Cells.Item(2284, 6).Value = 233.98
$ws.Cells.Item(2284, 7).Value = "Water and Natural Resources Committee and Staff"
$ws.Cells.Item(2284, 8).Value = 43474.5152412037
$ws.Cells.Item(2284, 9).Value = 2
$ws.Cells.Item(2284, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2284, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2284, 8).NumberFormat = "mm-dd-yy"

# Row 2285
$ws.Cells.Item(2285, 1).Value = 2485
$ws.Cells.Item(2285, 2).Value = 8550
$ws.Cells.Item(2285, 3).Value = "RSTP Dinner"
$ws.Cells.Item(2285, 4).Value = "Bull Ring, Santa Fe"
$ws.Cells.Item(2285, 5).Value = 43451
$ws.Cells.Item(2285, 6).Value = 235.49
$ws.Cells.Item(2285, 7).Value = "RSTP members and guests and staff"
$ws.Cells.Item(2285, 8).Value = 43475.8918245718
$ws.Cells.Item(2285, 9).Value = 2
$ws.Cells.Item(2285, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2285, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2285, 8).NumberFormat = "mm-dd-yy"

# Row 2286
$ws.Cells.Item(2286, 1).Value = 2486
$ws.Cells.Item(2286, 2).Value = 8558
$ws.Cells.Item(2286, 3).Value = "Rural Economic Development Committee"
$ws.Cells.Item(2286, 4).Value = "Rio Chama"
$ws.Cells.Item(2286, 5).Value = 43404
$ws.Cells.Item(2286, 6).Value = 221
$ws.Cells.Item(2286, 7).Value = "RuralEconomic Development Committee members and staff"
$ws.Cells.Item(2286, 8).Value = 43476
$ws.Cells.Item(2286, 9).Value = 2
$ws.Cells.Item(2286, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2286, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2286, 8).NumberFormat = "mm-dd-yy"

# Row 2287
$ws.Cells.Item(2287, 1).Value = 2488
$ws.Cells.Item(2287, 2).Value = 8589
$ws.Cells.Item(2287, 3).Value = "Breakfast"
$ws.Cells.Item(2287, 4).Value = "State Capitol - Room #307"
$ws.Cells.Item(2287, 5).Value = 43445
$ws.Cells.Item(2287, 6).Value = 341.58
$ws.Cells.Item(2287, 7).Value = "Legislative Finance Committee "
$ws.Cells.Item(2287, 8).Value = 43476
$ws.Cells.Item(2287, 9).Value = 2
$ws.Cells.Item(2287, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2287, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2287, 8).NumberFormat = "mm-dd-yy"

# Row 2288
$ws.Cells.Item(2288, 1).Value = 2489
$ws.Cells.Item(2288, 2).Value = 8589
$ws.Cells.Item(2288, 3).Value = "Breakfast"
$ws.Cells.Item(2288, 4).Value = "State Capitol - Room #322"
$ws.Cells.Item(2288, 5).Value = 43451
$ws.Cells.Item(2288, 6).Value = 374.11
$ws.Cells.Item(2288, 7).Value = "Revenue Stabilization Committee "
$ws.Cells.Item(2288, 8).Value = 43476.591280706
$ws.Cells.Item(2288, 9).Value = 2
$ws.Cells.Item(2288, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2288, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2288, 8).NumberFormat = "mm-dd-yy"

# Row 2289
$ws.Cells.Item(2289, 1).Value = 2490
$ws.Cells.Item(2289, 2).Value = 8609
$ws.Cells.Item(2289, 3).Value = "Governor's Ball Committee"
$ws.Cells.Item(2289, 4).Value = "Santa Fe, NM"
$ws.Cells.Item(2289, 5).Value = 43464
$ws.Cells.Item(2289, 6).Value = 1000
$ws.Cells.Item(2289, 7).Value = "Public"
$ws.Cells.Item(2289, 8).Value = 43477.5182101042
$ws.Cells.Item(2289, 9).Value = 2
$ws.Cells.Item(2289, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2289, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2289, 8).NumberFormat = "mm-dd-yy"

# Row 2290
$ws.Cells.Item(2290, 1).Value = 2491
$ws.Cells.Item(2290, 2).Value = 8656
$ws.Cells.Item(2290, 3).Value = "Legislative Committee Meeting"
$ws.Cells.Item(2290, 4).Value = "UNM Science and Technology Park Rotunda"
$ws.Cells.Item(2290, 5).Value = 43367
$ws.Cells.Item(2290, 6).Value = 1141.59
$ws.Cells.Item(2290, 7).Value = "Courts, Correction and Justice Committee, general public 9/24/2018 and 9/25/2018"
$ws.Cells.Item(2290, 8).Value = 43479.5166252662
$ws.Cells.Item(2290, 9).Value = 2
$ws.Cells.Item(2290, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2290, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2290, 8).NumberFormat = "mm-dd-yy"

# Row 2291
$ws.Cells.Item(2291, 1).Value = 2492
$ws.Cells.Item(2291, 2).Value = 8656
$ws.Cells.Item(2291, 3).Value = "Legislative Interim Committee Meeting"
$ws.Cells.Item(2291, 4).Value = "UNM Science and Technology Park Rotunda"
$ws.Cells.Item(2291, 5).Value = 43391
$ws.Cells.Item(2291, 6).Value = 1467.48
$ws.Cells.Item(2291, 7).Value = "Courts, Corrections and Justice Interim Committee, general public_x000d_`n"
$ws.Cells.Item(2291, 8).Value = 43479.5181403935
$ws.Cells.Item(2291, 9).Value = 2
$ws.Cells.Item(2291, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2291, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2291, 8).NumberFormat = "mm-dd-yy"

# Row 2292
$ws.Cells.Item(2292, 1).Value = 2493
$ws.Cells.Item(2292, 2).Value = 8656
$ws.Cells.Item(2292, 3).Value = "Interim Committee- Science Technology & Telecommun"
$ws.Cells.Item(2292, 4).Value = "Lobo Rain Forest"
$ws.Cells.Item(2292, 5).Value = 43391
$ws.Cells.Item(2292, 6).Value = 298
$ws.Cells.Item(2292, 7).Value = "interim committee, general public"
$ws.Cells.Item(2292, 8).Value = 43479.5188565972
$ws.Cells.Item(2292, 9).Value = 2
$ws.Cells.Item(2292, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2292, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2292, 8).NumberFormat = "mm-dd-yy"

# Row 2293
$ws.Cells.Item(2293, 1).Value = 2494
$ws.Cells.Item(2293, 2).Value = 8357
$ws.Cells.Item(2293, 3).Value = "House Freshmen Reception"
$ws.Cells.Item(2293, 4).Value = "Rio Chama"
$ws.Cells.Item(2293, 5).Value = 43440
$ws.Cells.Item(2293, 6).Value = 205.19
$ws.Cells.Item(2293, 7).Value = "Freshmen democrat legislators"
$ws.Cells.Item(2293, 8).Value = 43479.6928092245
$ws.Cells.Item(2293, 9).Value = 2
$ws.Cells.Item(2293, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2293, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2293, 8).NumberFormat = "mm-dd-yy"

# Row 2294
$ws.Cells.Item(2294, 1).Value = 2495
$ws.Cells.Item(2294, 2).Value = 8357
$ws.Cells.Item(2294, 3).Value = "Committee Dinner"
$ws.Cells.Item(2294, 4).Value = "Bull Ring"
$ws.Cells.Item(2294, 5).Value = 43433
$ws.Cells.Item(2294, 6).Value = 170
$ws.Cells.Item(2294, 7).Value = "Water and Natural Resources committee"
$ws.Cells.Item(2294, 8).Value = 43479
$ws.Cells.Item(2294, 9).Value = 2
$ws.Cells.Item(2294, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2294, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2294, 8).NumberFormat = "mm-dd-yy"

# Row 2295
$ws.Cells.Item(2295, 1).Value = 2496
$ws.Cells.Item(2295, 2).Value = 8357
$ws.Cells.Item(2295, 3).Value = "Committee Dinner"
$ws.Cells.Item(2295, 4).Value = "Bull Ring"
$ws.Cells.Item(2295, 5).Value = 43423
$ws.Cells.Item(2295, 6).Value = 75.53
$ws.Cells.Item(2295, 7).Value = "Military and Veterans affairs committee"
$ws.Cells.Item(2295, 8).Value = 43479.6963180556
$ws.Cells.Item(2295, 9).Value = 2
$ws.Cells.Item(2295, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2295, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2295, 8).NumberFormat = "mm-dd-yy"

# Row 2296
$ws.Cells.Item(2296, 1).Value = 2497
$ws.Cells.Item(2296, 2).Value = 8395
$ws.Cells.Item(2296, 3).Value = "Legislative Dinner"
$ws.Cells.Item(2296, 4).Value = "Rio Rancho"
$ws.Cells.Item(2296, 5).Value = 43438
$ws.Cells.Item(2296, 6).Value = 229.23
$ws.Cells.Item(2296, 7).Value = "Rio Rancho Public School Legislative Delegation _x000d_`nRep. Jason Harper_x000d_`nRep. Lewis_x000d_`nRep. Powdrell Culbert_x000d_`nRep. Ely_x000d_`nSen. Gould_x000d_`nSen. Brandt_x000d_`nSen. Sapien_x000d_`n"
$ws.Cells.Item(2296, 8).Value = 43480.4299287384
$ws.Cells.Item(2296, 9).Value = 2
$ws.Cells.Item(2296, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2296, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2296, 8).NumberFormat = "mm-dd-yy"

# Row 2297
$ws.Cells.Item(2297, 1).Value = 2498
$ws.Cells.Item(2297, 2).Value = 8730
$ws.Cells.Item(2297, 3).Value = "In Kind Event"
$ws.Cells.Item(2297, 4).Value = "Kelly Brew Pub"
$ws.Cells.Item(2297, 5).Value = 43383
$ws.Cells.Item(2297, 6).Value = 144.14
$ws.Cells.Item(2297, 7).Value = "Brian Egolf Fundraiser_x000d_`n"
$ws.Cells.Item(2297, 8).Value = 43480.443031331
$ws.Cells.Item(2297, 9).Value = 2
$ws.Cells.Item(2297, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2297, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2297, 8).NumberFormat = "mm-dd-yy"

# Row 2298
$ws.Cells.Item(2298, 1).Value = 2499
$ws.Cells.Item(2298, 2).Value = 8748
$ws.Cells.Item(2298, 3).Value = "Military and Veterans Affairs Interim Committee"
$ws.Cells.Item(2298, 4).Value = "Bull Ring"
$ws.Cells.Item(2298, 5).Value = 43425
$ws.Cells.Item(2298, 6).Value = 75.53
$ws.Cells.Item(2298, 7).Value = "Military and Veterans Affairs Interim Committee members and staff"
$ws.Cells.Item(2298, 8).Value = 43480.5546202199
$ws.Cells.Item(2298, 9).Value = 2
$ws.Cells.Item(2298, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2298, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2298, 8).NumberFormat = "mm-dd-yy"

# Row 2299
$ws.Cells.Item(2299, 1).Value = 2500
$ws.Cells.Item(2299, 2).Value = 8755
$ws.Cells.Item(2299, 3).Value = "Dinner"
$ws.Cells.Item(2299, 4).Value = "The Bull Ring"
$ws.Cells.Item(2299, 5).Value = 43451
$ws.Cells.Item(2299, 6).Value = 235.49
$ws.Cells.Item(2299, 7).Value = "Members and Staff of Revenue Stabilization & Tax Policy Committee"
$ws.Cells.Item(2299, 8).Value = 43480.6444996875
$ws.Cells.Item(2299, 9).Value = 2
$ws.Cells.Item(2299, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2299, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2299, 8).NumberFormat = "mm-dd-yy"

# Row 2300
$ws.Cells.Item(2300, 1).Value = 2501
$ws.Cells.Item(2300, 2).Value = 8682
$ws.Cells.Item(2300, 3).Value = "Lunch"
$ws.Cells.Item(2300, 4).Value = "La Chama "
$ws.Cells.Item(2300, 5).Value = 43398
$ws.Cells.Item(2300, 6).Value = 147.98
$ws.Cells.Item(2300, 7).Value = "Joint Finance Committee"
$ws.Cells.Item(2300, 8).Value = 43480.6787833333
$ws.Cells.Item(2300, 9).Value = 2
$ws.Cells.Item(2300, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2300, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2300, 8).NumberFormat = "mm-dd-yy"

# Row 2301
$ws.Cells.Item(2301, 1).Value = 2502
$ws.Cells.Item(2301, 2).Value = 8399
$ws.Cells.Item(2301, 3).Value = "Revenue Stabilization Tax Policy Committee Dinner"
$ws.Cells.Item(2301, 4).Value = "The Bull Ring"
$ws.Cells.Item(2301, 5).Value = 43424
$ws.Cells.Item(2301, 6).Value = 75.53
$ws.Cells.Item(2301, 7).Value = "Committee member and Staff"
$ws.Cells.Item(2301, 8).Value = 43480
$ws.Cells.Item(2301, 9).Value = 2
$ws.Cells.Item(2301, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2301, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2301, 8).NumberFormat = "mm-dd-yy"

# Row 2302
$ws.Cells.Item(2302, 1).Value = 2503
$ws.Cells.Item(2302, 2).Value = 8399
$ws.Cells.Item(2302, 3).Value = "Water and Natural Resources Dinner"
$ws.Cells.Item(2302, 4).Value = "The Bull Ring"
$ws.Cells.Item(2302, 5).Value = 43434
$ws.Cells.Item(2302, 6).Value = 250
$ws.Cells.Item(2302, 7).Value = "Members of the WNRC and Staff"
$ws.Cells.Item(2302, 8).Value = 43480.8840506944
$ws.Cells.Item(2302, 9).Value = 2
$ws.Cells.Item(2302, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2302, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2302, 8).NumberFormat = "mm-dd-yy"

# Row 2303
$ws.Cells.Item(2303, 1).Value = 2504
$ws.Cells.Item(2303, 2).Value = 8783
$ws.Cells.Item(2303, 3).Value = "Michael Padilla 7th Annual Matanza"
$ws.Cells.Item(2303, 4).Value = "7421 Isleta SW, Albuquerque"
$ws.Cells.Item(2303, 5).Value = 43386
$ws.Cells.Item(2303, 6).Value = 250
$ws.Cells.Item(2303, 7).Value = "Everyone"
$ws.Cells.Item(2303, 8).Value = 43480.9656375
$ws.Cells.Item(2303, 9).Value = 2
$ws.Cells.Item(2303, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2303, 6).NumberFormat = "#,##0.00"
$ws.Cells.Item(2303, 8).NumberFormat = "mm-dd-yy"

# --- Update defined name range to cover new extent ---
$n = $wb.Names.Item("Cam_SpecialEventLobbyist")
$n.RefersTo = "='Cam_SpecialEventLobbyist'!`$A`$1:`$I`$2303"
